$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.680.39'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.83%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.725.63'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9978'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.52'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9982'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4928'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.22%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.728.88'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '15.83'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06975'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6109'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.490'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.16'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9983'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.495.73'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9981'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007224'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.03%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.952.15'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.460'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.081'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '137.63'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.29'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.769'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.386'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.05'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.913'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.25%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.677'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04476'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.611'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.001'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6233'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9348'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.86%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.17%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9982'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01508'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.603'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.37'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3849'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.867'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1156'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05380'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.853'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.18'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '51.54'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.227'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.92%  '
